# Fix the truncated algorithm name in the results table: cell G3 ("name"
# column, row for the "TA Ciclica Aleatoria Ext" experiment) should read
# the full name "TA Ciclica Aleatoria Extendida".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "TA Ciclica Aleatoria Extendida"

# Reflect the cell the user was editing as the active selection.
$ws.Range("G3").Select() | Out-Null
